# Resume workbook update for HES application
# - Add "PySpark" to the Languages skill list
# - Add a new "Science in Africa" project (ProjectsDetails + ProjectsBullets)
# - Record the project year (2021) for the "Editors as Gatekeepers of Science" project

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. PersonalDetails: update the Languages line to mention PySpark
# ---------------------------------------------------------------------------
$wsPersonal = $wb.Worksheets.Item("PersonalDetails")
$wsPersonal.Range("B12").Value = "R, Python, SQL, LaTex, Stata, PySpark"

# ---------------------------------------------------------------------------
# 2. ProjectsDetails: fill in the missing project year for "Editors" project
#    and append a new row for the "Science in Africa" project
# ---------------------------------------------------------------------------
$wsProjects = $wb.Worksheets.Item("ProjectsDetails")
$wsProjects.Cells.Item(5, 6).Value = 2021

$wsProjects.Cells.Item(7, 1).Value = 6
$wsProjects.Cells.Item(7, 2).Value = "Science in Africa"
$wsProjects.Cells.Item(7, 3).Value = "AfricanScience"
$wsProjects.Cells.Item(7, 4).Value = "HBS"
$wsProjects.Cells.Item(7, 5).Value = "Harvard Business School"

# ---------------------------------------------------------------------------
# 3. ProjectsBullets: append the bullet describing the new project
# ---------------------------------------------------------------------------
$wsBullets = $wb.Worksheets.Item("ProjectsBullets")
$wsBullets.Cells.Item(13, 1).Value = 6
$wsBullets.Cells.Item(13, 2).Value = "AfricanScience"
$wsBullets.Cells.Item(13, 3).Value = "Pulled author information from a semi-strucutred API using Python"

# ---------------------------------------------------------------------------
# 4. Restore cursor/selection positions on each sheet touched during editing
#    (ProjectsBullets ends up being the active sheet/tab)
# ---------------------------------------------------------------------------
$wsPersonal.Range("N48").Select()

$wsWorkDeprecated = $wb.Worksheets.Item("WorkBulletsDepreciated")
$wsWorkDeprecated.Range("B12").Select()

$wsProjects.Range("C7").Select()

$wsBullets.Range("B14").Select()
